# Update "想去人数" (interest count) figures in column F across the
# workbook's sheets, reflecting refreshed crawl data as of commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 726
$ws.Range("F3").Value = 2781
$ws.Range("F4").Value = 1328
$ws.Range("F5").Value = 61
$ws.Range("F7").Value = 585
$ws.Range("F8").Value = 45
$ws.Range("F10").Value = 280
$ws.Range("F12").Value = 11572
$ws.Range("F13").Value = 6586
$ws.Range("F16").Value = 415
$ws.Range("F20").Value = 915
$ws.Range("F21").Value = 70
$ws.Range("F22").Value = 265
$ws.Range("F23").Value = 922
$ws.Range("F24").Value = 3637
$ws.Range("F26").Value = 984
$ws.Range("F27").Value = 495
$ws.Range("F28").Value = 165
$ws.Range("F29").Value = 312
$ws.Range("F30").Value = 18
$ws.Range("F32").Value = 296
$ws.Range("F33").Value = 5001
$ws.Range("F34").Value = 39
$ws.Range("F35").Value = 1231
$ws.Range("F36").Value = 228
$ws.Range("F37").Value = 428
$ws.Range("F38").Value = 190
$ws.Range("F39").Value = 532

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F14").Value = 10
$ws.Range("F24").Value = 38

# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 9032
$ws.Range("F3").Value = 496
$ws.Range("F4").Value = 1820

# Sheet 4: 全部类型 (All types)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 9032
$ws.Range("F3").Value = 496
$ws.Range("F4").Value = 1820
$ws.Range("F5").Value = 728
$ws.Range("F6").Value = 2781
$ws.Range("F10").Value = 1328
$ws.Range("F12").Value = 61
$ws.Range("F13").Value = 585
$ws.Range("F14").Value = 45
$ws.Range("F17").Value = 280
$ws.Range("F19").Value = 11572
$ws.Range("F24").Value = 415
$ws.Range("F28").Value = 915
$ws.Range("F29").Value = 70
$ws.Range("F30").Value = 265
$ws.Range("F31").Value = 922
$ws.Range("F32").Value = 3637
$ws.Range("F34").Value = 984
$ws.Range("F35").Value = 165
$ws.Range("F36").Value = 312
$ws.Range("F38").Value = 10
$ws.Range("F40").Value = 39
$ws.Range("F41").Value = 1231
$ws.Range("F42").Value = 228
$ws.Range("F43").Value = 190
$ws.Range("F44").Value = 532
$ws.Range("F48").Value = 38
